$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "분"
$ws.Range("C2").Value = 319
$ws.Range("B3").Value = "감사"
$ws.Range("C3").Value = 278
$ws.Range("B4").Value = "버스"
$ws.Range("C4").Value = 270
$ws.Range("B5").Value = "없"
$ws.Range("C5").Value = 251
$ws.Range("B6").Value = "보안"
$ws.Range("C6").Value = 213
$ws.Range("B7").Value = "좋"
$ws.Range("C7").Value = 199
$ws.Range("B8").Value = "글"
$ws.Range("C8").Value = 181
$ws.Range("B9").Value = "사람"
$ws.Range("C9").Value = 171
$ws.Range("B10").Value = "이용"
$ws.Range("C10").Value = 167
$ws.Range("B11").Value = "많"
$ws.Range("C11").Value = 144
$ws.Range("B12").Value = "생각"
$ws.Range("C12").Value = 144
$ws.Range("B13").Value = "직원"
$ws.Range("C13").Value = 144
$ws.Range("B14").Value = "문의"
$ws.Range("C14").Value = 138
$ws.Range("B15").Value = "관련"
$ws.Range("C15").Value = 122
$ws.Range("B16").Value = "사원"
$ws.Range("C16").Value = 115
$ws.Range("B17").Value = "회사"
$ws.Range("C17").Value = 110
$ws.Range("B18").Value = "사내"
$ws.Range("C18").Value = 108
$ws.Range("B19").Value = "차량"
$ws.Range("C19").Value = 108
$ws.Range("B20").Value = "하이닉스"
$ws.Range("C20").Value = 103
$ws.Range("B21").Value = "불편"
$ws.Range("C21").Value = 102
$ws.Range("B22").Value = "시"
$ws.Range("C22").Value = 102
$ws.Range("B23").Value = "시간"
$ws.Range("C23").Value = 101
$ws.Range("B24").Value = "사항"
$ws.Range("C24").Value = 97
$ws.Range("B25").Value = "옥계"
$ws.Range("C25").Value = 97
$ws.Range("B26").Value = "식당"
$ws.Range("C26").Value = 92
$ws.Range("B27").Value = "부탁"
$ws.Range("C27").Value = 89
$ws.Range("B28").Value = "사용"
$ws.Range("C28").Value = 89
$ws.Range("B29").Value = "답변"
$ws.Range("C29").Value = 89
$ws.Range("B30").Value = "확인"
$ws.Range("C30").Value = 89
$ws.Range("B31").Value = "일"
$ws.Range("C31").Value = 88
$ws.Range("B32").Value = "통근"
$ws.Range("C32").Value = 86
$ws.Range("B33").Value = "주차장"
$ws.Range("C33").Value = 85
$ws.Range("B34").Value = "말"
$ws.Range("C34").Value = 82
$ws.Range("B35").Value = "하"
$ws.Range("C35").Value = 82
$ws.Range("B36").Value = "문"
$ws.Range("C36").Value = 81
$ws.Range("B37").Value = "안녕"
$ws.Range("C37").Value = 80
$ws.Range("B38").Value = "안"
$ws.Range("C38").Value = 76
$ws.Range("B39").Value = "셔틀"
$ws.Range("C39").Value = 75
$ws.Range("B40").Value = "개선"
$ws.Range("C40").Value = 75
$ws.Range("B41").Value = "담당자"
$ws.Range("C41").Value = 75
$ws.Range("B42").Value = "요원"
$ws.Range("C42").Value = 74
$ws.Range("B43").Value = "분실물"
$ws.Range("C43").Value = 72
$ws.Range("B44").Value = "기숙사"
$ws.Range("C44").Value = 71
$ws.Range("B45").Value = "전"
$ws.Range("C45").Value = 70
$ws.Range("B46").Value = "근무"
$ws.Range("C46").Value = 70
$ws.Range("B47").Value = "출근"
$ws.Range("C47").Value = 70
$ws.Range("B48").Value = "오늘"
$ws.Range("C48").Value = 68
$ws.Range("B49").Value = "후"
$ws.Range("C49").Value = 68
$ws.Range("B50").Value = "앞"
$ws.Range("C50").Value = 65
$ws.Range("B51").Value = "이스텍"
$ws.Range("C51").Value = 65
$ws.Range("B52").Value = "안녕하"
$ws.Range("C52").Value = 64
$ws.Range("B53").Value = "친절"
$ws.Range("C53").Value = 63
$ws.Range("B54").Value = "치료"
$ws.Range("C54").Value = 63
$ws.Range("B55").Value = "요청"
$ws.Range("C55").Value = 63
$ws.Range("B56").Value = "중"
$ws.Range("C56").Value = 62
$ws.Range("B57").Value = "퇴근"
$ws.Range("C57").Value = 62
$ws.Range("B58").Value = "등"
$ws.Range("C58").Value = 62
$ws.Range("B59").Value = "어떻"
$ws.Range("C59").Value = 61
$ws.Range("B60").Value = "헬스장"
$ws.Range("C60").Value = 60
$ws.Range("B61").Value = "가능"
$ws.Range("C61").Value = 60
$ws.Range("B62").Value = "관리"
$ws.Range("C62").Value = 60
$ws.Range("B63").Value = "구성원"
$ws.Range("C63").Value = 59
$ws.Range("B64").Value = "칭찬"
$ws.Range("C64").Value = 58
$ws.Range("B65").Value = "문제"
$ws.Range("C65").Value = 57
$ws.Range("B66").Value = "화장실"
$ws.Range("C66").Value = 56
$ws.Range("B67").Value = "조치"
$ws.Range("C67").Value = 56
$ws.Range("B68").Value = "가방"
$ws.Range("C68").Value = 56
$ws.Range("B69").Value = "서비스"
$ws.Range("C69").Value = 56
$ws.Range("B70").Value = "분실"
$ws.Range("C70").Value = 55
$ws.Range("B71").Value = "하세"
$ws.Range("C71").Value = 55
$ws.Range("B72").Value = "흡연"
$ws.Range("C72").Value = 55
$ws.Range("B73").Value = "수고"
$ws.Range("C73").Value = 54
$ws.Range("B74").Value = "연락"
$ws.Range("C74").Value = 54
$ws.Range("B75").Value = "아이"
$ws.Range("C75").Value = 54
$ws.Range("B76").Value = "체육관"
$ws.Range("C76").Value = 54
$ws.Range("B77").Value = "청주"
$ws.Range("C77").Value = 53
$ws.Range("B78").Value = "검색"
$ws.Range("C78").Value = 53
$ws.Range("B79").Value = "정문"
$ws.Range("C79").Value = 53
$ws.Range("B80").Value = "안내"
$ws.Range("C80").Value = 53
$ws.Range("B81").Value = "교체"
$ws.Range("C81").Value = 52
$ws.Range("B82").Value = "운동"
$ws.Range("C82").Value = 52
$ws.Range("B83").Value = "부분"
$ws.Range("C83").Value = 52
$ws.Range("B84").Value = "모습"
$ws.Range("C84").Value = 51
$ws.Range("B85").Value = "남자"
$ws.Range("C85").Value = 51
$ws.Range("B86").Value = "업무"
$ws.Range("C86").Value = 51
$ws.Range("B87").Value = "기분"
$ws.Range("C87").Value = 50
$ws.Range("B88").Value = "전화"
$ws.Range("C88").Value = 50
$ws.Range("B89").Value = "경우"
$ws.Range("C89").Value = 50
$ws.Range("B90").Value = "예약"
$ws.Range("C90").Value = 50
$ws.Range("B91").Value = "시설"
$ws.Range("C91").Value = 49
$ws.Range("B92").Value = "고객"
$ws.Range("C92").Value = 49
$ws.Range("B93").Value = "좌석"
$ws.Range("C93").Value = 49
$ws.Range("B94").Value = "이렇"
$ws.Range("C94").Value = 49
$ws.Range("B95").Value = "검사"
$ws.Range("C95").Value = 49
$ws.Range("B96").Value = "건의"
$ws.Range("C96").Value = 49
$ws.Range("B97").Value = "곳"
$ws.Range("C97").Value = 48
$ws.Range("B98").Value = "터"
$ws.Range("C98").Value = 48
$ws.Range("B99").Value = "고생"
$ws.Range("C99").Value = 48
$ws.Range("B100").Value = "운영"
$ws.Range("C100").Value = 48
$ws.Range("B101").Value = "소리"
$ws.Range("C101").Value = 48
